$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (copy formatting from the existing header style)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New data values for columns I (I0) and J (IF)
$data = @(
    @(8, 9),
    @(1, 4),
    @(1, 6),
    @(1, 5),
    @(1, 6),
    @(1, 6),
    @(1, 6),
    @(1, 6),
    @(9, 9),
    @(9, 9),
    @(7, 8),
    @(8, 9),
    @(9, 9),
    @(9, 9),
    @(8, 8),
    @(6, 6),
    @(5, 6),
    @(1, 3),
    @(1, 4),
    @(7, 9),
    @(3, 4),
    @(3, 4),
    @(3, 3)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
